# Fruta / hortaliza, semanal
# The data rows (2-19) of the "Espárragos" sheet get re-shuffled: each
# target row receives the D,H,I,J,K,L,M,N,O,P values that used to belong
# to a different source row (row 13 is untouched). Columns A,B,C,E,F,G,Q,R
# are constant across all rows so they do not need to change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# target row -> source row (values copied from the source row "before" the edit)
$mapping = @{
    2  = 12
    3  = 9
    4  = 16
    5  = 17
    6  = 18
    7  = 14
    8  = 11
    9  = 8
    10 = 4
    11 = 5
    12 = 7
    13 = 13
    14 = 10
    15 = 2
    16 = 6
    17 = 3
    18 = 19
    19 = 15
}

$cols = @("D", "H", "I", "J", "K", "L", "M", "N", "O", "P")

# First snapshot the original values for every row/column we care about,
# so subsequent writes never read already-overwritten data.
$snapshot = @{}
foreach ($r in 2..19) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Now write back according to the mapping.
foreach ($targetRow in 2..19) {
    $sourceRow = $mapping[$targetRow]
    $sourceVals = $snapshot[$sourceRow]
    foreach ($col in $cols) {
        $ws.Range("$col$targetRow").Value = $sourceVals[$col]
    }
}
